$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update the greeting text for rule R10 (cell E8) as part of the git update.
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the new active cell/selection on the sheet.
$ws.Range("E8").Select()
